# chore: end work 18 and init work 19
# init work day 18, 11, 2022

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Estudos")

# --- Row 58 (18/11/2022): finish out the day — end time + rest period ---
$ws.Range("C58").Value = 0.84930555555555554
$ws.Range("E58").Value = 0.18680555555555556

# --- Row 59 (19/11/2022): start the next day + subject ---
$ws.Range("B59").Value = 0.71180555555555547
$ws.Range("H59").Value = "Atividade voluntária no Alpha EdTech"

# --- Update the ASSUNTO label for row 58 to include the new activities ---
$ws.Range("H58").Value = "Estágio + HARD + INGLÊS"

# --- Recalculate so DIF / ÚTEIS formula cells refresh ---
$excel.Calculate()

# --- Update the view state to match where the user ended up working ---
$excel.ActiveWindow.ScrollRow = 52
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("B60").Select()
